$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Reverser")
$rng = $ws.Range("B1:I1")
$rng.Interior.Color = 255
$rng.Font.Color = 16777215
Write-Host "done"
